# Atualização de bases das ligas, do dia: 28-05-2024 às 20:56
#
# The source feed re-ordered several match rows (the sequential "id"
# column in A stays put, but every other field for a pair of rows was
# swapped - including the database id in column B, the two team-name
# columns E/F, and all of the score / odds columns G:AD). Reproduce
# that by swapping the B:AD contents between each affected row pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Singapore Premier League")

$rowPairs = @(
    @(4, 5),
    @(18, 19),
    @(22, 23),
    @(26, 27),
    @(28, 29)
)

foreach ($pair in $rowPairs) {
    $row1 = $pair[0]
    $row2 = $pair[1]

    $range1 = $ws.Range("B" + $row1 + ":AD" + $row1)
    $range2 = $ws.Range("B" + $row2 + ":AD" + $row2)

    $values1 = $range1.Value()
    $values2 = $range2.Value()

    $range1.Value = $values2
    $range2.Value = $values1
}
